$d = $word.ActiveDocument

# 1. Update the delivery date from 01/12 to 06/07.
$d.Content.Find.Execute("01/12", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "06/07", 2)

# 2. Relocate the "_GoBack" bookmark from right after "...feito em grupos"
#    to the very end of the "Todos os trabalhos ... monitor da disciplina."
#    paragraph (Word re-anchors "_GoBack" to the location of the most
#    recent edit). We land there by inserting a throw-away marker
#    character right after that paragraph's final run, anchoring the new
#    bookmark just before it, and then removing the marker again -- this
#    sidesteps the runtime's quirky handling of zero-length ranges that
#    sit exactly on a paragraph's trailing boundary.
$target = $d.Content
$target.Find.Execute("monitor da disciplina.", $true, $false, $false, `
                      $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$anchorPos = $target.Start

$target.InsertAfter("~")

$bmRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$marker = $d.Range($anchorPos, $anchorPos + 1)
$marker.Delete()
